$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D18")
$sortKey = $ws.Range("A2:A18")

$dataRange.Sort($sortKey, 1)
